$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$ws.Range("G3").Value = "wait(3);`nvalidate1;`nlink_Click(network22_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0749_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;`nwifi_Mode(OFF);`nwait(2);`npress_Key(Home);`nlaunch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);`nvalidate5;`nwifi_Mode(ON);`nCheckUITextContains(Connected);`npress_Key(Home);"
$ws.Range("H3").Value = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=Network`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-0749`n};`nvalidate4`n{`nvalidate_Result=1`nvalidate_Result=Connected`n};`nvalidate5`n{`nvalidate_Result=2`nvalidate_Result=Disconnected`n};"
$ws.Range("G4").Value = "wait(3);`nvalidate1;`nlink_Click(network22_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0750_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;`nwifi_Mode(OFF);`nwait(2);`npress_Key(Home);`nlaunch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);`nvalidate5;`nwifi_Mode(ON);`nCheckUITextContains(Connected);`npress_Key(Home);"
$ws.Range("H4").Value = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=Network`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-0750`n};`nvalidate4`n{`nvalidate_doesNotContain=Connected`n};`nvalidate5`n{`nvalidate_doesNotContain=Connected`n};"
$ws.Range("G5").Value = "wait(3);`nvalidate1;`nlink_Click(network22_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0751_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;`nwifi_Mode(OFF);`nwait(2);`npress_Key(Home);`nlaunch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);`nvalidate5;`nwifi_Mode(ON);`nCheckUITextContains(Connected);`npress_Key(Home);"
$ws.Range("H5").Value = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=Network`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-0751`n};`nvalidate4`n{`nvalidate_Result=1`nvalidate_Result=Connected`n};`nvalidate5`n{`nvalidate_Result=2`nvalidate_Result=Disconnected`n};"
$ws.Range("G6").Value = "wait(3);`nvalidate1;`nlink_Click(network22_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0752_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;`nwifi_Mode(OFF);`nwait(2);`npress_Key(Home);`nlaunch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);`nvalidate5;`nwifi_Mode(ON);`nCheckUITextContains(Connected);`npress_Key(Home);"
$ws.Range("H6").Value = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=Network`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-0752`n};`nvalidate4`n{`nvalidate_Result=Connected`n};`nvalidate5`n{`nvalidate_Result=Disconnected`n};"
$ws.Range("G7").Value = "wait(3);`nvalidate1;`nlink_Click(network22_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0755_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;`nwifi_Mode(OFF);`nwait(60);`npress_Key(Home);`nlaunch_App_Device(com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);`nvalidate5;`nwifi_Mode(ON);`nCheckUITextContains(Connected);"
$ws.Range("H7").Value = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=Network`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-0755`n};`nvalidate4`n{`nvalidate_Result=Connected`n};`nvalidate5`n{`nvalidate_Result=Disconnected`n};"
$ws.Range("G8").Value = "wait(3);`nvalidate1;`nlink_Click(network22_test_link);`nwait(2);`nvalidate2;`nSelectTestToRun(VT187_0757_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nvalidate4;`npress_Key(Home);"
$ws.Range("H8").Value = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=Network`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-0757`n};`nvalidate4`n{`nvalidate_PageTitle=Navigation Check`n};"

$null = $ws.Range("E1").Select()
